# Update 1-Apr-2021, end of day.
# Applies the 31-Mar-2021 / 1-Apr-2021 petty cash entries to "Sheet1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 17 (31/Mar/2021, "Wages Expense"): add a second disbursement ---
$ws.Range("D17").Formula = "=60000+7010000"

# --- Row 18 ("A/R"): add two more credits ---
$ws.Range("C18").Formula = "=100000000+18000000+38573000"

# --- Row 19 ("TRANSFER BCA"): add two more debits ---
$ws.Range("D19").Formula = "=100000000+2150000+28384000+609000"

# --- Row 20: new "SALES - cash/retail" entry ---
$ws.Range("B20").Value = "SALES - cash/retail"
$ws.Range("C20").Formula = "=12961475+31166525-38573000"

# --- Row 21: new "PRIVE - andreas" entry ---
$ws.Range("B21").Value = "PRIVE - andreas"
$ws.Range("D21").Value = 2000000

# --- Row 22: new "SELISIH - lebih" entry ---
$ws.Range("B22").Value = "SELISIH - lebih"
$ws.Range("C22").Value = 425000

# --- Row 23: new "PRIVE - bulanan" entry ---
$ws.Range("B23").Value = "PRIVE - bulanan"
$ws.Range("D23").Value = 16000000

# --- Row 24: new "SETOR KE BANK" entry ---
$ws.Range("B24").Value = "SETOR KE BANK"
$ws.Range("D24").Value = 6000000

# --- Row 25: new date, start of 1/Apr/2021 ---
$ws.Range("A25").Value = 44287

# --- Update view state to match where the user ended up ---
$ws.Range("C45").Select()
